$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.358.32"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$ws.Range("D3").Value = "1.842.44"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("D5").Value = "240.07"
$ws.Range("E5").Value = "  -0.28%  "

# Row 6
$ws.Range("D6").Value = "0.6307"
$ws.Range("E6").Value = "  +0.39%  "

# Row 7
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").Value = "0.07474"
$ws.Range("E8").Value = "  +0.45%  "

# Row 9
$ws.Range("D9").Value = "0.2900"
$ws.Range("E9").Value = "  +0.40%  "

# Row 10
$ws.Range("D10").Value = "24.97"
$ws.Range("E10").Value = "  +3.01%  "

# Row 11
$ws.Range("D11").Value = "0.07728"
$ws.Range("E11").Value = "  +0.03%  "

# Row 12
$ws.Range("D12").Value = "1.842.80"
$ws.Range("E12").Value = "  +0.03%  "

# Row 13
$ws.Range("D13").Value = "4.981"
$ws.Range("E13").Value = "  +0.03%  "

# Row 14
$ws.Range("E14").Value = "  +0.24%  "

# Row 15
$ws.Range("D15").Value = "0.00001031"
$ws.Range("E15").Value = "  +2.70%  "

# Row 16
$ws.Range("D16").Value = "81.94"
$ws.Range("E16").Value = "  -0.07%  "

# Row 17
$ws.Range("D17").Value = "6.227"
$ws.Range("E17").Value = "  +1.66%  "

# Row 18
$ws.Range("D18").Value = "29.352.38"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
$ws.Range("D19").Value = "229.30"
$ws.Range("E19").Value = "  +0.83%  "

# Row 20
$ws.Range("D20").Value = "12.33"
$ws.Range("E20").Value = "  +0.61%  "

# Row 21
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").Value = "  -0.24%  "

# Row 22
$ws.Range("D22").Value = "7.395"
$ws.Range("E22").Value = "  +0.41%  "

# Row 23
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").Value = "  -0.30%  "

# Row 24
$ws.Range("D24").Value = "158.17"
$ws.Range("E24").Value = "  -0.28%  "

# Row 25
$ws.Range("D25").Value = "8.538"
$ws.Range("E25").Value = "  +1.92%  "

# Row 26
$ws.Range("E26").Value = "  -1.08%  "

# Row 27
$ws.Range("D27").Value = "17.48"
$ws.Range("E27").Value = "  -0.26%  "

# Row 28
$ws.Range("D28").Value = "0.06850"
$ws.Range("E28").Value = "  +9.14%  "

# Row 29
$ws.Range("E29").Value = "  +4.78%  "

# Row 30
$ws.Range("D30").Value = "1.489"
$ws.Range("E30").Value = "  +1.00%  "

# Row 31
$ws.Range("D31").Value = "4.075"
$ws.Range("E31").Value = "  +0.89%  "

# Row 32
$ws.Range("E32").Value = "  -0.30%  "

# Row 33
$ws.Range("D33").Value = "1.831"
$ws.Range("E33").Value = "  +0.84%  "

# Row 35
$ws.Range("D35").Value = "0.7005"
$ws.Range("E35").Value = "  +0.99%  "

# Row 36
$ws.Range("D36").Value = "2.585"
$ws.Range("E36").Value = "  -0.11%  "

# Row 37
$ws.Range("D37").Value = "0.01846"
$ws.Range("E37").Value = "  +1.92%  "

# Row 38
$ws.Range("D38").Value = "2.819"
$ws.Range("E38").Value = "  -0.26%  "

# Row 39
$ws.Range("D39").Value = "1.237.57"
$ws.Range("E39").Value = "  -0.96%  "

# Row 40
$ws.Range("D40").Value = "6.801"
$ws.Range("E40").Value = "  +4.45%  "

# Row 41
$ws.Range("D41").Value = "0.9444"
$ws.Range("E41").Value = "  +3.90%  "

# Row 42
$ws.Range("D42").Value = "0.9988"
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").Value = "1.997.50"
$ws.Range("E43").Value = "  -0.16%  "

# Row 44
$ws.Range("D44").Value = "100.99"
$ws.Range("E44").Value = "  -0.32%  "

# Row 45
$ws.Range("D45").Value = "65.51"
$ws.Range("E45").Value = "  -0.75%  "

# Row 46
$ws.Range("E46").Value = "  +3.56%  "

# Row 47
$ws.Range("D47").Value = "1.719"
$ws.Range("E47").Value = "  +3.83%  "

# Row 48
$ws.Range("D48").Value = "7.043"
$ws.Range("E48").Value = "  +0.09%  "

# Row 49
$ws.Range("D49").Value = "8.971"
$ws.Range("E49").Value = "  -0.41%  "

# Row 50
$ws.Range("E50").Value = "  -1.24%  "

# Row 51
$ws.Range("D51").Value = "0.3917"
$ws.Range("E51").Value = "  -0.32%  "
